$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before current row 15 (uv -> voter row), shifting rows 15+ down by 1.
$ws.Rows.Item(15).Insert()
$ws.Range("A15").Value = "ink"
$ws.Range("B15").Value = "uv"
$ws.Range("C15").Value = "ink->uv"
$ws.Range("D15").Value = "shows under"
$ws.Range("E15").Value = 1

# Insert two new rows after the (now shifted) "uv -> voter" row, which is row 16.
$ws.Rows.Item(17).Insert()
$ws.Range("A17").Value = "elections"
$ws.Range("B17").Value = "significance"
$ws.Range("C17").Value = "elections->significance"
$ws.Range("D17").Value = "are assuming"
$ws.Range("E17").Value = 1

$ws.Rows.Item(18).Insert()
$ws.Range("A18").Value = "elections"
$ws.Range("B18").Value = "prelude"
$ws.Range("C18").Value = "elections->prelude"
$ws.Range("D18").Value = "are"
$ws.Range("E18").Value = 1

# Append a new row at the end (row 35, after the old last row 31 shifted to 34).
$ws.Range("A35").Value = "ink"
$ws.Range("B35").Value = "finger"
$ws.Range("C35").Value = "ink->finger"
$ws.Range("D35").Value = "stays on"
$ws.Range("E35").Value = 1
